# "Support seq and : function" -- add a new Nvidia/Intel/ATI CPU-vs-GPU
# comparison table (rows 13-16) to Sheet4, then leave Sheet2 as the
# active/selected sheet instead of Sheet4.

$wb  = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item(4)

# Enter the new comparison block. The order below matters: it reproduces
# the exact order new strings were first introduced into the shared
# string table (Nvidia, Intel, previous, now, Optimal).
$ws4.Range("B15").Value = "Nvidia"
$ws4.Range("B14").Value = "Intel"
$ws4.Range("C13").Value = "previous"
$ws4.Range("D13").Value = "now"
$ws4.Range("E13").Value = "Optimal"
$ws4.Range("F13").Value = "CPU"
$ws4.Range("B16").Value = "ATI"

$ws4.Range("C14").Value = 9.92
$ws4.Range("D14").Value = 1.71
$ws4.Range("E14").Value = 4.55
$ws4.Range("F14").Value = 3.295

$ws4.Range("C15").Value = 2.8
$ws4.Range("D15").Value = 1.85
$ws4.Range("E15").Value = 1.78

$ws4.Range("C16").Value = 5
$ws4.Range("D16").Value = 3.56

# Sheet4's selection moves off the old N10 cell onto the new table.
[void]$ws4.Range("C25").Select()

# The workbook's active tab moves from Sheet4 (index 3) to Sheet2 (index 1).
$ws2 = $wb.Worksheets.Item(2)
[void]$ws2.Activate()
